$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(452, 1).Value2 = "Login with valid username and password"
$ws.Cells.Item(452, 2).Value2 = "PASSED"
$ws.Cells.Item(452, 3).Value2 = "chrome"

$ws.Cells.Item(453, 1).Value2 = "Login with valid username and password"
$ws.Cells.Item(453, 2).Value2 = "PASSED"
$ws.Cells.Item(453, 3).Value2 = "chrome"

$ws.Cells.Item(454, 1).Value2 = "Login with valid username and password"
$ws.Cells.Item(454, 2).Value2 = "PASSED"
$ws.Cells.Item(454, 3).Value2 = "chrome"

$ws.Cells.Item(455, 1).Value2 = "Login with valid username and password"
$ws.Cells.Item(455, 2).Value2 = "PASSED"
$ws.Cells.Item(455, 3).Value2 = "chrome"

$ws.Cells.Item(456, 1).Value2 = "Login with valid username and password"
$ws.Cells.Item(456, 2).Value2 = "PASSED"
$ws.Cells.Item(456, 3).Value2 = "chrome"

$ws.Cells.Item(457, 1).Value2 = "Login with valid username and password"
$ws.Cells.Item(457, 2).Value2 = "PASSED"
$ws.Cells.Item(457, 3).Value2 = "chrome"

$ws.Cells.Item(458, 1).Value2 = "Login with valid username and password"
$ws.Cells.Item(458, 2).Value2 = "PASSED"
$ws.Cells.Item(458, 3).Value2 = "chrome"

$ws.Cells.Item(459, 1).Value2 = "Create a country"
$ws.Cells.Item(459, 2).Value2 = "PASSED"
$ws.Cells.Item(459, 3).Value2 = "chrome"

$ws.Cells.Item(460, 1).Value2 = "Create and Delete Cities"
$ws.Cells.Item(460, 2).Value2 = "FAILED"
$ws.Cells.Item(460, 3).Value2 = "chrome"

$ws.Cells.Item(461, 1).Value2 = "Login with valid username and password"
$ws.Cells.Item(461, 2).Value2 = "PASSED"
$ws.Cells.Item(461, 3).Value2 = "chrome"

$ws.Cells.Item(462, 1).Value2 = "Create a country"
$ws.Cells.Item(462, 2).Value2 = "PASSED"
$ws.Cells.Item(462, 3).Value2 = "chrome"

$ws.Cells.Item(463, 1).Value2 = "Create a country"
$ws.Cells.Item(463, 2).Value2 = "PASSED"
$ws.Cells.Item(463, 3).Value2 = "chrome"

$ws.Cells.Item(464, 1).Value2 = "Create a country 2"
$ws.Cells.Item(464, 2).Value2 = "PASSED"
$ws.Cells.Item(464, 3).Value2 = "chrome"

$ws.Cells.Item(465, 1).Value2 = "Create a citizenship"
$ws.Cells.Item(465, 2).Value2 = "PASSED"
$ws.Cells.Item(465, 3).Value2 = "chrome"

$ws.Cells.Item(466, 1).Value2 = "Create a Citizenship"
$ws.Cells.Item(466, 2).Value2 = "PASSED"
$ws.Cells.Item(466, 3).Value2 = "chrome"

$ws.Cells.Item(467, 1).Value2 = "Create a Citizenship"
$ws.Cells.Item(467, 2).Value2 = "PASSED"
$ws.Cells.Item(467, 3).Value2 = "chrome"

$ws.Cells.Item(468, 1).Value2 = "Create a Citizenship"
$ws.Cells.Item(468, 2).Value2 = "PASSED"
$ws.Cells.Item(468, 3).Value2 = "chrome"

$ws.Cells.Item(469, 1).Value2 = "Users List"
$ws.Cells.Item(469, 2).Value2 = "PASSED"

$ws.Cells.Item(470, 1).Value2 = "Create Country"
$ws.Cells.Item(470, 2).Value2 = "PASSED"
$ws.Cells.Item(470, 3).Value2 = "chrome"

$ws.Cells.Item(471, 1).Value2 = "Create Nationality"
$ws.Cells.Item(471, 2).Value2 = "PASSED"
$ws.Cells.Item(471, 3).Value2 = "chrome"

$ws.Cells.Item(472, 1).Value2 = "Fee Functionality"
$ws.Cells.Item(472, 2).Value2 = "PASSED"
$ws.Cells.Item(472, 3).Value2 = "chrome"

$ws.Cells.Item(473, 1).Value2 = "Fee Functionality"
$ws.Cells.Item(473, 2).Value2 = "PASSED"
$ws.Cells.Item(473, 3).Value2 = "chrome"

$ws.Cells.Item(474, 1).Value2 = "Create Nationality and Delete"
$ws.Cells.Item(474, 2).Value2 = "PASSED"
$ws.Cells.Item(474, 3).Value2 = "chrome"

$ws.Cells.Item(475, 1).Value2 = "Create Nationality and Delete"
$ws.Cells.Item(475, 2).Value2 = "PASSED"
$ws.Cells.Item(475, 3).Value2 = "chrome"

$ws.Cells.Item(476, 1).Value2 = "Create and Delete Cities"
$ws.Cells.Item(476, 2).Value2 = "FAILED"
$ws.Cells.Item(476, 3).Value2 = "chrome"

$ws.Cells.Item(477, 1).Value2 = "Create and Delete Cities"
$ws.Cells.Item(477, 2).Value2 = "FAILED"
$ws.Cells.Item(477, 3).Value2 = "chrome"

$ws.Cells.Item(478, 1).Value2 = "Create and Delete Cities"
$ws.Cells.Item(478, 2).Value2 = "FAILED"
$ws.Cells.Item(478, 3).Value2 = "chrome"

$ws.Cells.Item(479, 1).Value2 = "Create and Delete Cities"
$ws.Cells.Item(479, 2).Value2 = "FAILED"
$ws.Cells.Item(479, 3).Value2 = "chrome"

$ws.Cells.Item(480, 1).Value2 = "Create CitizenShip and Delete"
$ws.Cells.Item(480, 2).Value2 = "PASSED"
$ws.Cells.Item(480, 3).Value2 = "chrome"

$ws.Cells.Item(481, 1).Value2 = "Create Inventory and Delete"
$ws.Cells.Item(481, 2).Value2 = "FAILED"
$ws.Cells.Item(481, 3).Value2 = "chrome"

$ws.Cells.Item(482, 1).Value2 = "Create Inventory and Delete"
$ws.Cells.Item(482, 2).Value2 = "FAILED"
$ws.Cells.Item(482, 3).Value2 = "chrome"

$ws.Cells.Item(483, 1).Value2 = "Create and Delete CitizenShip From Excel"
$ws.Cells.Item(483, 2).Value2 = "FAILED"
$ws.Cells.Item(483, 3).Value2 = "chrome"

$ws.Cells.Item(484, 1).Value2 = "States testing with JDBC"
$ws.Cells.Item(484, 2).Value2 = "FAILED"
$ws.Cells.Item(484, 3).Value2 = "chrome"

$ws.Cells.Item(485, 1).Value2 = "States testing with JDBC"
$ws.Cells.Item(485, 2).Value2 = "FAILED"
$ws.Cells.Item(485, 3).Value2 = "chrome"

$ws.Cells.Item(486, 1).Value2 = "Create and Delete Cities"
$ws.Cells.Item(486, 2).Value2 = "FAILED"
$ws.Cells.Item(486, 3).Value2 = "chrome"

$ws.Cells.Item(487, 1).Value2 = "Create and Delete Cities"
$ws.Cells.Item(487, 2).Value2 = "FAILED"
$ws.Cells.Item(487, 3).Value2 = "chrome"

$ws.Cells.Item(488, 1).Value2 = "Fees create and delete functionality"
$ws.Cells.Item(488, 2).Value2 = "PASSED"
$ws.Cells.Item(488, 3).Value2 = "chrome"

$ws.Cells.Item(489, 1).Value2 = "Create a Country"
$ws.Cells.Item(489, 2).Value2 = "PASSED"
$ws.Cells.Item(489, 3).Value2 = "chrome"

$ws.Cells.Item(490, 1).Value2 = "Create a CitizenShip"
$ws.Cells.Item(490, 2).Value2 = "PASSED"
$ws.Cells.Item(490, 3).Value2 = "chrome"

$ws.Cells.Item(491, 1).Value2 = "Create Inventory and Delete"
$ws.Cells.Item(491, 2).Value2 = "FAILED"
$ws.Cells.Item(491, 3).Value2 = "chrome"

$ws.Cells.Item(492, 1).Value2 = "Create Inventory and Delete"
$ws.Cells.Item(492, 2).Value2 = "FAILED"
$ws.Cells.Item(492, 3).Value2 = "chrome"

$ws.Cells.Item(493, 1).Value2 = "Create Inventory and Delete"
$ws.Cells.Item(493, 2).Value2 = "FAILED"
$ws.Cells.Item(493, 3).Value2 = "chrome"

$ws.Cells.Item(494, 1).Value2 = "Create Inventory and Delete"
$ws.Cells.Item(494, 2).Value2 = "FAILED"
$ws.Cells.Item(494, 3).Value2 = "chrome"

$ws.Cells.Item(495, 1).Value2 = "Create Inventory and Delete"
$ws.Cells.Item(495, 2).Value2 = "FAILED"
$ws.Cells.Item(495, 3).Value2 = "chrome"

$ws.Cells.Item(496, 1).Value2 = "Create Inventory and Delete"
$ws.Cells.Item(496, 2).Value2 = "FAILED"
$ws.Cells.Item(496, 3).Value2 = "chrome"

$ws.Cells.Item(497, 1).Value2 = "Create Inventory and Delete"
$ws.Cells.Item(497, 2).Value2 = "FAILED"
$ws.Cells.Item(497, 3).Value2 = "chrome"

$ws.Cells.Item(498, 1).Value2 = "Create Inventory and Delete"
$ws.Cells.Item(498, 2).Value2 = "FAILED"
$ws.Cells.Item(498, 3).Value2 = "chrome"

$ws.Cells.Item(499, 1).Value2 = "Create Inventory and Delete"
$ws.Cells.Item(499, 2).Value2 = "FAILED"
$ws.Cells.Item(499, 3).Value2 = "chrome"

$ws.Cells.Item(500, 1).Value2 = "Create Inventory and Delete"
$ws.Cells.Item(500, 2).Value2 = "FAILED"
$ws.Cells.Item(500, 3).Value2 = "chrome"

$ws.Cells.Item(501, 1).Value2 = "Create Inventory and Delete"
$ws.Cells.Item(501, 2).Value2 = "FAILED"
$ws.Cells.Item(501, 3).Value2 = "chrome"

$ws.Cells.Item(502, 1).Value2 = "Create Inventory and Delete"
$ws.Cells.Item(502, 2).Value2 = "FAILED"
$ws.Cells.Item(502, 3).Value2 = "chrome"

$ws.Cells.Item(503, 1).Value2 = "Create Inventory and Delete"
$ws.Cells.Item(503, 2).Value2 = "FAILED"
$ws.Cells.Item(503, 3).Value2 = "chrome"

$ws.Cells.Item(504, 1).Value2 = "Create Inventory and Delete"
$ws.Cells.Item(504, 2).Value2 = "PASSED"
$ws.Cells.Item(504, 3).Value2 = "chrome"
